$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 23 with the TSC in_ex values to the cell_data table
$ws.Range("A23").Value = 2527.0553673936838
$ws.Range("B23").Value = 1675.9277600756905
$ws.Range("C23").Value = 1776.6481060779727

# Keep the selection in sync with the new data extent
$ws.Range("A1:C23").Select()
